$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. Excel shifts everything
# (including formulas) one column to the right, e.g. the old
# "=D15+B12" formula in D6 automatically becomes "=E15+C12" in E6.
$ws.Range("A1").EntireColumn.Insert()

# Carry the header formatting (bold style) from the former A1 (now B1)
# onto the brand-new A1 cell.
$ws.Range("B1").Copy($ws.Range("A1"))

# New headers: "ID" (row index) in column A, "Group" (old column A
# header) stays as the header text, now living in column B.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Group"

# Fill in the new "ID" column (just the record/row number) and the new
# "Group" column values for each of the 8 data rows.
$ids = @(1, 2, 3, 4, 5, 6, 7, 8)
$groups = @(1, 1, 1, 2, 2, 1, 2, 1)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = $groups[$i]
}

# Match the saved selection/active cell.
$ws.Range("B10").Select() | Out-Null
